# The "Gender" column (F) for the 4th data row (row 5) was re-labeled from
# the informal "Girl" to "Female" in the shared-strings table. Apply that
# as a plain cell-value edit so the workbook's shared-strings table and the
# cell's <v> index get rebuilt accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "Female"

# The author's cursor ended up on G7 when the edit was committed/saved.
$null = $ws.Range("G7").Select()
